# Fruta / hortaliza, semanal
# Insert two new weekly observation rows (559 and 560) into the
# "Hortaliza, Vega Central Mapocho de Santiago - Zapallo italiano" data sheet,
# pushing the previously existing rows 559-582 down to 561-584.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 559 (shifts old 559:582 -> 561:584)
$ws.Rows("559:560").Insert()

# ---- New row 559 ----
$ws.Cells.Item(559,1).Value  = 9
$ws.Cells.Item(559,2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(559,3).Value  = "Metropolitana"
$ws.Cells.Item(559,4).Value  = 45075
$ws.Cells.Item(559,5).Value  = 13
$ws.Cells.Item(559,6).Value  = 100112032
$ws.Cells.Item(559,7).Value  = "Zapallo italiano"
$ws.Cells.Item(559,8).Value  = "Bola 8"
$ws.Cells.Item(559,9).Value  = "Primera"
$ws.Cells.Item(559,10).Value = 52
$ws.Cells.Item(559,11).Value = 10000
$ws.Cells.Item(559,12).Value = 11000
$ws.Cells.Item(559,13).Value = 10500
$ws.Cells.Item(559,14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(559,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(559,16).Value = 210
$ws.Cells.Item(559,17).Value = 50
$ws.Cells.Item(559,18).Value = "Hortaliza"

# ---- New row 560 ----
$ws.Cells.Item(560,1).Value  = 9
$ws.Cells.Item(560,2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(560,3).Value  = "Metropolitana"
$ws.Cells.Item(560,4).Value  = 45075
$ws.Cells.Item(560,5).Value  = 13
$ws.Cells.Item(560,6).Value  = 100112032
$ws.Cells.Item(560,7).Value  = "Zapallo italiano"
$ws.Cells.Item(560,8).Value  = "Sin especificar"
$ws.Cells.Item(560,9).Value  = "Primera"
$ws.Cells.Item(560,10).Value = 97
$ws.Cells.Item(560,11).Value = 10000
$ws.Cells.Item(560,12).Value = 12000
$ws.Cells.Item(560,13).Value = 10969
$ws.Cells.Item(560,14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(560,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(560,16).Value = 219
$ws.Cells.Item(560,17).Value = 50
$ws.Cells.Item(560,18).Value = "Hortaliza"
